$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 102: add a "No" marker in column E (existing shared string) ---
$ws.Range("E102").Value = "No"

# --- New data rows 109-111, written in an order that reproduces the ---
# --- original shared-string append order from the source workbook.  ---

# Row 109
$ws.Range("A109").Value = "3/5/2025``"
$ws.Range("B109").Value = "Amazon"
$ws.Range("C109").Value = "sr. data scientist"
$ws.Range("D109").Value = "Interview w/ management"

# Row 110
$ws.Range("B110").Value = "DTN"
$ws.Range("C110").Value = "Sr Data Scientist"
$ws.Range("D110").Value = "First interview"

# Row 111
$ws.Range("B111").Value = "MISO"
$ws.Range("C111").Value = "Senior Engineer Market Evaluation"
$ws.Range("D111").Value = "application, MISO market efficiency evaluator"

# URLs / trailing notes, filled back-to-front (111 -> 110 -> 109)
$ws.Range("F111").Value = "https://recruiting.ultipro.com/MID1029MISO/JobBoard/362b6b1d-f1c3-46f5-9554-4aa90e2bda64/OpportunityDetail?opportunityId=ae933cb8-99a0-42af-ae01-b41431abde3a"
$ws.Range("F110").Value = "same url as before"
$ws.Range("F109").Value = "the energy job, same url as before"

# Remaining cells referencing already-existing shared strings
$ws.Range("E109").Value = "No"

# Dates for rows 110/111 (numeric serials, formatted like the rest of column A)
$ws.Range("A110").Value = 45723
$ws.Range("A110").NumberFormat = "m/d/yy"
$ws.Range("A111").Value = 45724
$ws.Range("A111").NumberFormat = "m/d/yy"

$ws.Range("F109").Select()
